$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from G1 (an existing "sum" header) onto H1, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New "Save" column values for the two data rows.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
